$wb = $excel.ActiveWorkbook

# 1. Update B8 on the "current file" sheet: "Egypt" -> "Egypt, Cairo"
$cur = $wb.Worksheets.Item("current file")
$cur.Range("B8").Value = "Egypt, Cairo"
$cur.Range("B10").Select()

# 2. Duplicate "Sheet1" -> creates "Sheet1 (2)" right after "Sheet1", before "Poposed file"
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Copy([System.Reflection.Missing]::Value, $sheet1)

# 3. On the new copy, tweak cell B5 (date) and selection
$copy = $wb.Worksheets.Item("Sheet1 (2)")
$copy.Range("B5").Value = 43515
$copy.Range("H8").Select()

# 4. On "Poposed file" sheet, change the frozen-pane top-left cell and selection
$prop = $wb.Worksheets.Item("Poposed file")
$prop.Select()
$prop.Range("Y15").Select()
